$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows for "RM 232" (row 26) and "SC 92" (row 28).
# Delete the later row first so the earlier row index stays valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# After the deletions, apply the remaining per-cell value updates (mostly in
# column F, a couple in column C) to match the target state. Cells that go
# "missing" need to stay present as empty text cells (matching the existing
# empty-result placeholders elsewhere in the sheet) rather than being cleared
# outright, so a lone "'" (Excel's force-text prefix) is used to write an
# empty text cell instead of assigning "" (which removes the cell entirely).
# The style is reset back to Normal afterwards so the quote-prefix formatting
# that "'" implies doesn't stick around.
function Set-Missing($row, $col) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'"
    $c.Style = "Normal"
}

$ws.Cells.Item(2, 6).Value = 18.03
Set-Missing 6 6
$ws.Cells.Item(12, 6).Value = 17.45
Set-Missing 14 6
$ws.Cells.Item(20, 6).Value = 17.73
$ws.Cells.Item(21, 6).Value = 16.58
Set-Missing 23 6
Set-Missing 24 6
$ws.Cells.Item(26, 3).Value = 10.8
Set-Missing 27 3
$ws.Cells.Item(30, 3).Value = 11.4
$ws.Cells.Item(31, 6).Value = 17.18
Set-Missing 32 3
$ws.Cells.Item(33, 6).Value = 17.53
